$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("FEINmismatch")
$ws1.Range("B2").Value = "Tue Nov 12 17:20:00 EST 2024"
$ws1.Range("B3").Value = "Tue Nov 12 17:20:14 EST 2024"
$ws1.Range("B4").Value = "Tue Nov 12 17:20:26 EST 2024"
$ws1.Range("B5").Value = "Tue Nov 12 17:20:40 EST 2024"
$ws1.Range("B6").Value = "Tue Nov 12 17:20:52 EST 2024"
$ws1.Range("B7").Value = "Tue Nov 12 17:21:04 EST 2024"
$ws1.Range("B8").Value = "Tue Nov 12 17:21:15 EST 2024"
$ws1.Range("B9").Value = "Tue Nov 12 17:21:26 EST 2024"
$ws1.Range("B10").Value = "Tue Nov 12 17:21:38 EST 2024"
$ws1.Range("B11").Value = "Tue Nov 12 17:21:50 EST 2024"
$ws1.Range("B12").Value = "Tue Nov 12 17:22:01 EST 2024"
$ws1.Range("B13").Value = "Tue Nov 12 17:22:13 EST 2024"
$ws1.Range("B14").Value = "Tue Nov 12 17:22:27 EST 2024"
$ws1.Range("B15").Value = "Tue Nov 12 17:22:39 EST 2024"
$ws1.Range("B16").Value = "Tue Nov 12 17:22:50 EST 2024"
$ws1.Range("B17").Value = "Tue Nov 12 17:23:01 EST 2024"
$ws1.Range("B18").Value = "Tue Nov 12 17:23:12 EST 2024"
$ws1.Range("B19").Value = "Tue Nov 12 17:23:23 EST 2024"
$ws1.Range("B20").Value = "Tue Nov 12 17:23:34 EST 2024"
$ws1.Range("B21").Value = "Tue Nov 12 17:23:48 EST 2024"
$ws1.Range("B22").Value = "Tue Nov 12 17:23:59 EST 2024"
$ws1.Range("B23").Value = "Tue Nov 12 17:24:10 EST 2024"
$ws1.Range("B24").Value = "Tue Nov 12 17:24:22 EST 2024"
$ws1.Range("B25").Value = "Tue Nov 12 17:24:33 EST 2024"
$ws1.Range("B26").Value = "Tue Nov 12 17:24:48 EST 2024"
$ws1.Range("B27").Value = "Tue Nov 12 17:24:59 EST 2024"
$ws1.Range("B28").Value = "Tue Nov 12 17:25:10 EST 2024"
$ws1.Range("B29").Value = "Tue Nov 12 17:25:21 EST 2024"
$ws1.Range("B30").Value = "Tue Nov 12 17:25:33 EST 2024"

$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")
$ws2.Range("B2").Value = "Tue Nov 12 17:25:49 EST 2024"
$ws2.Range("B3").Value = "Tue Nov 12 17:26:00 EST 2024"
$ws2.Range("B4").Value = "Tue Nov 12 17:26:11 EST 2024"
$ws2.Range("B5").Value = "Tue Nov 12 17:26:22 EST 2024"
$ws2.Range("B6").Value = "Tue Nov 12 17:26:32 EST 2024"
$ws2.Range("B7").Value = "Tue Nov 12 17:26:47 EST 2024"
$ws2.Range("B8").Value = "Tue Nov 12 17:26:58 EST 2024"
$ws2.Range("B9").Value = "Tue Nov 12 17:27:09 EST 2024"
$ws2.Range("B10").Value = "Tue Nov 12 17:27:19 EST 2024"
$ws2.Range("B11").Value = "Tue Nov 12 17:27:30 EST 2024"
$ws2.Range("B12").Value = "Tue Nov 12 17:27:42 EST 2024"
$ws2.Range("B13").Value = "Tue Nov 12 17:27:53 EST 2024"
$ws2.Range("B14").Value = "Tue Nov 12 17:28:04 EST 2024"
$ws2.Range("B15").Value = "Tue Nov 12 17:28:15 EST 2024"
$ws2.Range("B16").Value = "Tue Nov 12 17:28:25 EST 2024"
$ws2.Range("B17").Value = "Tue Nov 12 17:28:37 EST 2024"
$ws2.Range("B18").Value = "Tue Nov 12 17:28:49 EST 2024"
$ws2.Range("B19").Value = "Tue Nov 12 17:29:00 EST 2024"
